# Regenerate s_val data to filter save games: update B2:E7 and G2:G7 values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    3 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    4 = @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 2.652525447291612)
    5 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    6 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    7 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
